$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
SELECT DISTINCT
    sts.study_subject_id AS "Case ID",
    prg.program_acronym AS "Program Code",
    prg.program_id AS "Program ID",
    std.study_acronym AS "ARM",
    sts.disease_subtype AS "Diagnosis",
    stf.grouped_recurrence_score AS "Recurrence Score",
    dgn.tumor_size_group AS "Tumor Size (cm)",
    dgn.er_status AS "ER Status",
    dgn.pr_status AS "PR Status",
    dem.age_at_index AS "Age (years)",
    CASE 
        WHEN dem.survival_time = 0 THEN ''
        ELSE dem.survival_time
    END AS "Survival (days)"
FROM 
    df_program prg
LEFT JOIN 
    df_study std ON prg.program_id = std."program.program_id"
LEFT JOIN 
    df_study_subject sts ON std.study_id = sts."study.study_id"
LEFT JOIN
    df_diagnosis dgn ON sts.study_subject_id = dgn."study_subject.study_subject_id"
LEFT JOIN
    df_stratification_factor stf ON sts.study_subject_id = stf."study_subject.study_subject_id"
LEFT JOIN
    df_demographic_data dem ON sts.study_subject_id = dem."study_subject.study_subject_id"
WHERE 
    std.study_acronym = 'A'
ORDER BY 
    sts.study_subject_id ASC
LIMIT 100;
'@

$filesQuery = @'
SELECT DISTINCT
    fil.file_name AS "File Name",
    smp.type AS "Association",
    fil.file_description AS "Description",
    fil.file_format AS "File Format",
    CASE 
        WHEN ROUND(fil.file_size, 2) = CAST(ROUND(fil.file_size, 0) AS INT) 
        THEN CAST(CAST(ROUND(fil.file_size, 0) AS INT) AS TEXT) || ' Bytes'
        ELSE ROUND(fil.file_size, 2) || ' Bytes'
    END AS "Size",
    prg.program_acronym AS "Program Code",
    std.study_acronym AS "ARM",
    sts.study_subject_id AS "Case ID",
    smp.sample_id AS "Sample ID"
FROM 
    df_program prg
LEFT JOIN 
    df_study std ON prg.program_id = std."program.program_id"
LEFT JOIN 
    df_study_subject sts ON std.study_id = sts."study.study_id"
LEFT JOIN
    df_diagnosis dgn ON sts.study_subject_id = dgn."study_subject.study_subject_id"
LEFT JOIN
    df_stratification_factor stf ON sts.study_subject_id = stf."study_subject.study_subject_id"
LEFT JOIN
    df_demographic_data dem ON sts.study_subject_id = dem."study_subject.study_subject_id"
LEFT JOIN
    df_sample smp ON sts.study_subject_id = smp."study_subject.study_subject_id"
LEFT JOIN
    df_laboratory_procedure lab ON prg.program_id = lab."program.program_id"
LEFT JOIN
    df_file fil ON smp.sample_id = fil."sample.sample_id"
WHERE 
    std.study_acronym = 'A'
ORDER BY 
    fil.file_name ASC
LIMIT 100;
'@

# Update the Files tab query (B4) so file sizes render as a friendly "x Bytes"
# string instead of a raw (sometimes fractional) number.
$ws.Range("B4").Value = $filesQuery

# Update the Cases tab query (B2) to null-out zero survival times and drop the
# now-unused sample/lab/file joins.
$ws.Range("B2").Value = $casesQuery

# Leave the cursor on the cell that was actually edited, matching the saved
# selection state.
$ws.Range("B2").Select()

Write-Output "done"
